$d = $word.ActiveDocument

# --- 1. Insert a new "Search Books" row into the Backlog table, just
#        above the existing "Total" row. ------------------------------
$t = $d.Tables.Item(1)
$totalRow = $t.Rows.Item($t.Rows.Count)          # the "Total" row
$newRow = $t.Rows.Add($totalRow)                  # insert row above it

$newRow.Cells.Item(1).Range.Text = "1"
$newRow.Cells.Item(2).Range.Text = "Search Books"
$newRow.Cells.Item(3).Range.Text = "2"
$newRow.Cells.Item(4).Range.Text = "1"

# NOTE: after Rows.Add the old $totalRow handle now refers to the row
# at that same (now stale) position, i.e. the row we just inserted -
# re-fetch the real "Total" row (now the last row) fresh.
$totalRow = $t.Rows.Item($t.Rows.Count)

# --- 2. Update the Total (Size / Story points) cell 28 -> 29. ---------
#        Only the last digit actually changes (8 -> 9); editing just
#        that character keeps the existing "2" run untouched and
#        produces two adjacent runs ("2" and "9") instead of merging
#        them back into a single "29" run.
$totalCell = $totalRow.Cells.Item($totalRow.Cells.Count)
$totalRng = $totalCell.Range.Duplicate
$totalRng.Find.ClearFormatting()
$totalRng.Find.Execute("28") | Out-Null

$digitRng = $totalRng.Duplicate
$digitRng.Collapse(0) | Out-Null            # wdCollapseEnd
$digitRng.MoveStart(1, -1) | Out-Null       # select just the trailing "8"
$digitRng.Text = "9"

# Toggling formatting on the freshly typed digit stops the engine from
# silently re-merging it with the preceding "2" run (both end up with
# identical run properties again, but as two distinct <w:r> elements).
$digitRng.Font.Bold = 1
$digitRng.Font.Bold = 0

# --- 3. The extra row pushes the document onto a second page, so the
#        footer's cached PAGE field result has to move from 1 to 2. ----
$footer = $d.Sections.Item(1).Footers.Item(1)
$pageField = $footer.Range.Fields.Item(1)
$pageField.Result.Characters.Item(1).Text = "2"
